$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 12.81844533333333
$ws.Range("H2").Value = 38.455336
$ws.Range("I2").Value = 0.1651242884662064
$ws.Range("J2").Value = 0.1651242884662064
$ws.Range("M2").Value = 150.538648
$ws.Range("N2").Value = 451.615944
$ws.Range("O2").Value = 0.5643166737150654
$ws.Range("P2").Value = 0.5643166737150654
$ws.Range("Q2").Value = 1929.67142994191
$ws.Range("R2").Value = 17367.04286947719
$ws.Range("S2").Value = 0.09318238921681656
$ws.Range("T2").Value = 0.09318238921681651
$ws.Range("G3").Value = 12.81844533333333
$ws.Range("H3").Value = 38.455336
$ws.Range("I3").Value = 0.1651242884662064
$ws.Range("J3").Value = 0.1651242884662064
$ws.Range("O3").Value = 0.1124385517418854
$ws.Range("P3").Value = 0.1124385517418854
$ws.Range("Q3").Value = 384.4817476187378
$ws.Range("R3").Value = 3460.33572856864
$ws.Range("S3").Value = 0.01856633585254956
$ws.Range("T3").Value = 0.01856633585254956
$ws.Range("G4").Value = 12.81844533333333
$ws.Range("H4").Value = 38.455336
$ws.Range("I4").Value = 0.1651242884662064
$ws.Range("J4").Value = 0.1651242884662064
$ws.Range("M4").Value = 33.76251433333334
$ws.Range("N4").Value = 101.287543
$ws.Range("O4").Value = 0.1265638428268858
$ws.Range("P4").Value = 0.1265638428268858
$ws.Range("Q4").Value = 432.7829442977165
$ws.Range("R4").Value = 3895.046498679448
$ws.Range("S4").Value = 0.0208987644923383
$ws.Range("T4").Value = 0.02089876449233829
$ws.Range("G5").Value = 12.81844533333333
$ws.Range("H5").Value = 38.455336
$ws.Range("I5").Value = 0.1651242884662064
$ws.Range("J5").Value = 0.1651242884662064
$ws.Range("M5").Value = 17.946055
$ws.Range("N5").Value = 53.838165
$ws.Range("O5").Value = 0.06727347560546457
$ws.Range("P5").Value = 0.06727347560546457
$ws.Range("Q5").Value = 230.0405249664933
$ws.Range("R5").Value = 2070.36472469844
$ws.Range("S5").Value = 0.01110848479200103
$ws.Range("T5").Value = 0.01110848479200103
$ws.Range("G6").Value = 12.81844533333333
$ws.Range("H6").Value = 38.455336
$ws.Range("I6").Value = 0.1651242884662064
$ws.Range("J6").Value = 0.1651242884662064
$ws.Range("M6").Value = 34.52108433333333
$ws.Range("N6").Value = 103.563253
$ws.Range("O6").Value = 0.1294074561106987
$ws.Range("P6").Value = 0.1294074561106987
$ws.Range("Q6").Value = 442.5066323742232
$ws.Range("R6").Value = 3982.559691368008
$ws.Range("S6").Value = 0.02136831411250096
$ws.Range("T6").Value = 0.02136831411250095
$ws.Range("I7").Value = 0.2089670898510501
$ws.Range("J7").Value = 0.2089670898510501
$ws.Range("M7").Value = 150.538648
$ws.Range("N7").Value = 451.615944
$ws.Range("O7").Value = 0.5643166737150654
$ws.Range("P7").Value = 0.5643166737150654
$ws.Range("Q7").Value = 2442.026105482357
$ws.Range("R7").Value = 21978.23494934121
$ws.Range("S7").Value = 0.1179236130606618
$ws.Range("T7").Value = 0.1179236130606618
$ws.Range("I8").Value = 0.2089670898510501
$ws.Range("J8").Value = 0.2089670898510501
$ws.Range("O8").Value = 0.1124385517418854
$ws.Range("P8").Value = 0.1124385517418854
$ws.Range("Q8").Value = 486.5670135328177
$ws.Range("S8").Value = 0.02349595694456851
$ws.Range("T8").Value = 0.02349595694456851
$ws.Range("I9").Value = 0.2089670898510501
$ws.Range("J9").Value = 0.2089670898510501
$ws.Range("M9").Value = 33.76251433333334
$ws.Range("N9").Value = 101.287543
$ws.Range("O9").Value = 0.1265638428268858
$ws.Range("P9").Value = 0.1265638428268858
$ws.Range("Q9").Value = 547.6928515308724
$ws.Range("R9").Value = 4929.235663777851
$ws.Range("S9").Value = 0.02644767791590003
$ws.Range("T9").Value = 0.02644767791590002
$ws.Range("I10").Value = 0.2089670898510501
$ws.Range("J10").Value = 0.2089670898510501
$ws.Range("M10").Value = 17.946055
$ws.Range("N10").Value = 53.838165
$ws.Range("O10").Value = 0.06727347560546457
$ws.Range("P10").Value = 0.06727347560546457
$ws.Range("Q10").Value = 291.1194924536733
$ws.Range("R10").Value = 2620.07543208306
$ws.Range("S10").Value = 0.01405794242143954
$ws.Range("T10").Value = 0.01405794242143954
$ws.Range("I11").Value = 0.2089670898510501
$ws.Range("J11").Value = 0.2089670898510501
$ws.Range("M11").Value = 34.52108433333333
$ws.Range("N11").Value = 103.563253
$ws.Range("O11").Value = 0.1294074561106987
$ws.Range("P11").Value = 0.1294074561106987
$ws.Range("Q11").Value = 559.998314396699
$ws.Range("R11").Value = 5039.984829570291
$ws.Range("S11").Value = 0.0270418995084802
$ws.Range("T11").Value = 0.02704189950848019
$ws.Range("G12").Value = 29.03561066666667
$ws.Range("H12").Value = 87.106832
$ws.Range("I12").Value = 0.3740301126102598
$ws.Range("J12").Value = 0.3740301126102597
$ws.Range("M12").Value = 150.538648
$ws.Range("N12").Value = 451.615944
$ws.Range("O12").Value = 0.5643166737150654
$ws.Range("P12").Value = 0.5643166737150654
$ws.Range("Q12").Value = 4370.981573614378
$ws.Range("R12").Value = 39338.83416252941
$ws.Range("S12").Value = 0.2110714290174932
$ws.Range("T12").Value = 0.2110714290174931
$ws.Range("G13").Value = 29.03561066666667
$ws.Range("H13").Value = 87.106832
$ws.Range("I13").Value = 0.3740301126102598
$ws.Range("J13").Value = 0.3740301126102597
$ws.Range("O13").Value = 0.1124385517418854
$ws.Range("P13").Value = 0.1124385517418854
$ws.Range("Q13").Value = 870.9061077217422
$ws.Range("R13").Value = 7838.154969495679
$ws.Range("S13").Value = 0.04205540416975192
$ws.Range("T13").Value = 0.0420554041697519
$ws.Range("G14").Value = 29.03561066666667
$ws.Range("H14").Value = 87.106832
$ws.Range("I14").Value = 0.3740301126102598
$ws.Range("J14").Value = 0.3740301126102597
$ws.Range("M14").Value = 33.76251433333334
$ws.Range("N14").Value = 101.287543
$ws.Range("O14").Value = 0.1265638428268858
$ws.Range("P14").Value = 0.1265638428268858
$ws.Range("Q14").Value = 980.3152213104196
$ws.Range("R14").Value = 8822.836991793776
$ws.Range("S14").Value = 0.04733868838492733
$ws.Range("T14").Value = 0.0473386883849273
$ws.Range("G15").Value = 29.03561066666667
$ws.Range("H15").Value = 87.106832
$ws.Range("I15").Value = 0.3740301126102598
$ws.Range("J15").Value = 0.3740301126102597
$ws.Range("M15").Value = 17.946055
$ws.Range("N15").Value = 53.838165
$ws.Range("O15").Value = 0.06727347560546457
$ws.Range("P15").Value = 0.06727347560546457
$ws.Range("Q15").Value = 521.0746659825866
$ws.Range("R15").Value = 4689.67199384328
$ws.Range("S15").Value = 0.02516230565639548
$ws.Range("T15").Value = 0.02516230565639547
$ws.Range("G16").Value = 29.03561066666667
$ws.Range("H16").Value = 87.106832
$ws.Range("I16").Value = 0.3740301126102598
$ws.Range("J16").Value = 0.3740301126102597
$ws.Range("M16").Value = 34.52108433333333
$ws.Range("N16").Value = 103.563253
$ws.Range("O16").Value = 0.1294074561106987
$ws.Range("P16").Value = 0.1294074561106987
$ws.Range("Q16").Value = 1002.340764493833
$ws.Range("R16").Value = 9021.066880444496
$ws.Range("S16").Value = 0.0484022853816919
$ws.Range("T16").Value = 0.04840228538169188
$ws.Range("G17").Value = 3.868294000000001
$ws.Range("H17").Value = 11.604882
$ws.Range("I17").Value = 0.04983048082025044
$ws.Range("J17").Value = 0.04983048082025043
$ws.Range("M17").Value = 150.538648
$ws.Range("N17").Value = 451.615944
$ws.Range("O17").Value = 0.5643166737150654
$ws.Range("P17").Value = 0.5643166737150654
$ws.Range("Q17").Value = 582.3277488265121
$ws.Range("R17").Value = 5240.949739438609
$ws.Range("S17").Value = 0.0281201711861061
$ws.Range("T17").Value = 0.02812017118610609
$ws.Range("G18").Value = 3.868294000000001
$ws.Range("H18").Value = 11.604882
$ws.Range("I18").Value = 0.04983048082025044
$ws.Range("J18").Value = 0.04983048082025043
$ws.Range("O18").Value = 0.1124385517418854
$ws.Range("P18").Value = 0.1124385517418854
$ws.Range("Q18").Value = 116.0272091308533
$ws.Range("R18").Value = 1044.24488217768
$ws.Range("S18").Value = 0.005602867096030757
$ws.Range("T18").Value = 0.005602867096030756
$ws.Range("G19").Value = 3.868294000000001
$ws.Range("H19").Value = 11.604882
$ws.Range("I19").Value = 0.04983048082025044
$ws.Range("J19").Value = 0.04983048082025043
$ws.Range("M19").Value = 33.76251433333334
$ws.Range("N19").Value = 101.287543
$ws.Range("O19").Value = 0.1265638428268858
$ws.Range("P19").Value = 0.1265638428268858
$ws.Range("Q19").Value = 130.6033316205474
$ws.Range("R19").Value = 1175.429984584926
$ws.Range("S19").Value = 0.006306737142522325
$ws.Range("T19").Value = 0.006306737142522323
$ws.Range("G20").Value = 3.868294000000001
$ws.Range("H20").Value = 11.604882
$ws.Range("I20").Value = 0.04983048082025044
$ws.Range("J20").Value = 0.04983048082025043
$ws.Range("M20").Value = 17.946055
$ws.Range("N20").Value = 53.838165
$ws.Range("O20").Value = 0.06727347560546457
$ws.Range("P20").Value = 0.06727347560546457
$ws.Range("Q20").Value = 69.42061688017
$ws.Range("R20").Value = 624.7855519215301
$ws.Range("S20").Value = 0.003352269635869688
$ws.Range("T20").Value = 0.003352269635869687
$ws.Range("G21").Value = 3.868294000000001
$ws.Range("H21").Value = 11.604882
$ws.Range("I21").Value = 0.04983048082025044
$ws.Range("J21").Value = 0.04983048082025043
$ws.Range("M21").Value = 34.52108433333333
$ws.Range("N21").Value = 103.563253
$ws.Range("O21").Value = 0.1294074561106987
$ws.Range("P21").Value = 0.1294074561106987
$ws.Range("Q21").Value = 133.5377034001274
$ws.Range("R21").Value = 1201.839330601146
$ws.Range("S21").Value = 0.006448435759721574
$ws.Range("T21").Value = 0.006448435759721572
$ws.Range("G22").Value = 15.684801
$ws.Range("H22").Value = 47.05440299999999
$ws.Range("I22").Value = 0.2020480282522334
$ws.Range("J22").Value = 0.2020480282522333
$ws.Range("M22").Value = 150.538648
$ws.Range("N22").Value = 451.615944
$ws.Range("O22").Value = 0.5643166737150654
$ws.Range("P22").Value = 0.5643166737150654
$ws.Range("Q22").Value = 2361.168736689048
$ws.Range("R22").Value = 21250.51863020143
$ws.Range("S22").Value = 0.1140190712339879
$ws.Range("T22").Value = 0.1140190712339879
$ws.Range("G23").Value = 15.684801
$ws.Range("H23").Value = 47.05440299999999
$ws.Range("I23").Value = 0.2020480282522334
$ws.Range("J23").Value = 0.2020480282522333
$ws.Range("O23").Value = 0.1124385517418854
$ws.Range("P23").Value = 0.1124385517418854
$ws.Range("Q23").Value = 470.4564042450799
$ws.Range("R23").Value = 4234.107638205719
$ws.Range("S23").Value = 0.02271798767898466
$ws.Range("T23").Value = 0.02271798767898466
$ws.Range("G24").Value = 15.684801
$ws.Range("H24").Value = 47.05440299999999
$ws.Range("I24").Value = 0.2020480282522334
$ws.Range("J24").Value = 0.2020480282522333
$ws.Range("M24").Value = 33.76251433333334
$ws.Range("N24").Value = 101.287543
$ws.Range("O24").Value = 0.1265638428268858
$ws.Range("P24").Value = 0.1265638428268858
$ws.Range("Q24").Value = 529.558318577981
$ws.Range("R24").Value = 4766.024867201828
$ws.Range("S24").Value = 0.02557197489119785
$ws.Range("T24").Value = 0.02557197489119784
$ws.Range("G25").Value = 15.684801
$ws.Range("H25").Value = 47.05440299999999
$ws.Range("I25").Value = 0.2020480282522334
$ws.Range("J25").Value = 0.2020480282522333
$ws.Range("M25").Value = 17.946055
$ws.Range("N25").Value = 53.838165
$ws.Range("O25").Value = 0.06727347560546457
$ws.Range("P25").Value = 0.06727347560546457
$ws.Range("Q25").Value = 281.4803014100549
$ws.Range("R25").Value = 2533.322712690494
$ws.Range("S25").Value = 0.01359247309975884
$ws.Range("T25").Value = 0.01359247309975883
$ws.Range("G26").Value = 15.684801
$ws.Range("H26").Value = 47.05440299999999
$ws.Range("I26").Value = 0.2020480282522334
$ws.Range("J26").Value = 0.2020480282522333
$ws.Range("M26").Value = 34.52108433333333
$ws.Range("N26").Value = 103.563253
$ws.Range("O26").Value = 0.1294074561106987
$ws.Range("P26").Value = 0.1294074561106987
$ws.Range("Q26").Value = 541.4563380725509
$ws.Range("R26").Value = 4873.107042652959
$ws.Range("S26").Value = 0.02614652134830411
$ws.Range("T26").Value = 0.02614652134830409
